$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin price ("D") and volume-change ("E") values for this run.
# D is $null for rows where only the volume percentage changed.
$updates = @(
    [pscustomobject]@{ Row = 2; D = "27.766.81"; E = "  +1.00%  " },
    [pscustomobject]@{ Row = 3; D = "1.775.30"; E = "  +1.53%  " },
    [pscustomobject]@{ Row = 4; D = "1.002"; E = "  -0.07%  " },
    [pscustomobject]@{ Row = 5; D = "327.21"; E = "  +0.79%  " },
    [pscustomobject]@{ Row = 6; D = $null; E = "  +0.09%  " },
    [pscustomobject]@{ Row = 7; D = "0.4580"; E = "  +2.43%  " },
    [pscustomobject]@{ Row = 8; D = "0.3581"; E = "  -0.41%  " },
    [pscustomobject]@{ Row = 9; D = "0.07487"; E = "  -0.10%  " },
    [pscustomobject]@{ Row = 10; D = $null; E = "  -0.40%  " },
    [pscustomobject]@{ Row = 11; D = "1.104"; E = "  +1.02%  " },
    [pscustomobject]@{ Row = 12; D = $null; E = "  +0.01%  " },
    [pscustomobject]@{ Row = 13; D = "20.83"; E = "  +0.79%  " },
    [pscustomobject]@{ Row = 14; D = "6.039"; E = "  +0.30%  " },
    [pscustomobject]@{ Row = 15; D = "7.206"; E = "  +1.01%  " },
    [pscustomobject]@{ Row = 16; D = "1.776.58"; E = "  +1.59%  " },
    [pscustomobject]@{ Row = 17; D = "93.57"; E = "  +0.01%  " },
    [pscustomobject]@{ Row = 18; D = "0.00001060"; E = "  -0.10%  " },
    [pscustomobject]@{ Row = 19; D = "0.06437"; E = "  +0.85%  " },
    [pscustomobject]@{ Row = 20; D = "1.001"; E = "  +0.14%  " },
    [pscustomobject]@{ Row = 21; D = "17.05"; E = "  +1.57%  " },
    [pscustomobject]@{ Row = 22; D = $null; E = "  -0.82%  " },
    [pscustomobject]@{ Row = 23; D = "27.796.00"; E = "  +0.89%  " },
    [pscustomobject]@{ Row = 24; D = "11.32"; E = "  +0.72%  " },
    [pscustomobject]@{ Row = 25; D = $null; E = "  -0.04%  " },
    [pscustomobject]@{ Row = 26; D = "163.85"; E = "  +0.83%  " },
    [pscustomobject]@{ Row = 27; D = "20.27"; E = "  -0.90%  " },
    [pscustomobject]@{ Row = 28; D = "1.977.47"; E = "  +1.38%  " },
    [pscustomobject]@{ Row = 29; D = "2.178"; E = "  +4.54%  " },
    [pscustomobject]@{ Row = 30; D = "125.71"; E = "  +0.11%  " },
    [pscustomobject]@{ Row = 31; D = $null; E = "  +1.82%  " },
    [pscustomobject]@{ Row = 32; D = "0.09199"; E = "  +1.35%  " },
    [pscustomobject]@{ Row = 33; D = "3.671"; E = "  +0.17%  " },
    [pscustomobject]@{ Row = 34; D = "5.536"; E = "  -0.04%  " },
    [pscustomobject]@{ Row = 35; D = "11.86"; E = "  -0.80%  " },
    [pscustomobject]@{ Row = 36; D = "0.02291"; E = "  -0.16%  " },
    [pscustomobject]@{ Row = 37; D = "0.06169"; E = "  +2.55%  " },
    [pscustomobject]@{ Row = 38; D = "0.2091"; E = "  +0.27%  " },
    [pscustomobject]@{ Row = 39; D = "0.6323"; E = "  -0.71%  " },
    [pscustomobject]@{ Row = 40; D = "4.953"; E = "  +0.24%  " },
    [pscustomobject]@{ Row = 41; D = "1.185"; E = "  -1.56%  " },
    [pscustomobject]@{ Row = 42; D = "1.391"; E = "  +0.60%  " },
    [pscustomobject]@{ Row = 43; D = "7.809"; E = "  +0.44%  " },
    [pscustomobject]@{ Row = 44; D = "13.27"; E = "  +0.26%  " },
    [pscustomobject]@{ Row = 45; D = "3.740"; E = "  +0.43%  " },
    [pscustomobject]@{ Row = 46; D = "0.5919"; E = "  +0.37%  " },
    [pscustomobject]@{ Row = 47; D = "122.75"; E = "  +0.07%  " },
    [pscustomobject]@{ Row = 48; D = $null; E = "  -0.15%  " },
    [pscustomobject]@{ Row = 49; D = $null; E = "  +0.93%  " },
    [pscustomobject]@{ Row = 50; D = $null; E = "  -1.05%  " },
    [pscustomobject]@{ Row = 51; D = "72.34"; E = "  +0.18%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # The Price column holds plain text (e.g. "27.766.81", "1.002"), not numbers.
        # Force text format first so Excel does not auto-convert numeric-looking
        # strings into numbers (which would silently drop formatting such as
        # trailing zeros or turn very small numbers into scientific notation).
        $cellD = $ws.Range("D$($u.Row)")
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
